# feat: add whitelist-based logging for bot commands
#
# Inserts two new sub-activity rows ("API manajemen whitelist" and
# "UI React untuk manajemen whitelist") under the "Keamanan" stage right
# after row 18, shifting the following rows down by two. Marks the
# "Selesai" status (column F) of the Keamanan/Logging/Version-Control rows
# that are now complete (rows 17,18,21,22) as TRUE, which also applies the
# green "done" fill used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows after row 18 (pushes rows 19+ down to 21+)
$ws.Rows.Item(19).Resize(2).Insert()

# Copy the border/format from the row above onto the two new rows so they
# pick up the same per-cell borders used throughout the table (A:G).
$ws.Range("A18:G18").Copy()
$ws.Range("A19:A20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 19 - new sub-activity: API manajemen whitelist
$ws.Range("A19").Value = 17.1
$ws.Range("B19").Value = "Keamanan"
$ws.Range("C19").Value = "API manajemen whitelist"
$ws.Range("D19").Value = "CRUD whitelist"
$ws.Range("E19").Value = "Endpoint API"
$ws.Range("F19").Value = $true
$ws.Range("G19").Value = ""

# Row 20 - new sub-activity: UI React untuk manajemen whitelist
$ws.Range("A20").Value = 17.2
$ws.Range("B20").Value = "Keamanan"
$ws.Range("C20").Value = "UI React untuk manajemen whitelist"
$ws.Range("D20").Value = "UI whitelist"
$ws.Range("E20").Value = "Screenshot UI"
$ws.Range("F20").Value = $true
$ws.Range("G20").Value = ""

# Mark completed rows (now shifted): 17, 18 stay in place; former 19 & 20
# (Logging / Version Control) are now rows 21 & 22 after the insert.
$ws.Range("F17").Value = $true
$ws.Range("F18").Value = $true
$ws.Range("F21").Value = $true
$ws.Range("F22").Value = $true

# Apply the "done" green fill (same style used for other TRUE rows) to the
# rows that just flipped to TRUE, matching cells A:G.
$doneFill = $ws.Range("A2:G2")  # sample of an existing "done" (green) row
$doneRanges = @("A17:G17", "A18:G18", "A19:G20", "A21:G21", "A22:G22")
foreach ($rng in $doneRanges) {
    $ws.Range($rng).Interior.Color = $doneFill.Interior.Color
}

# Refresh the autofilter / used range to include the two new rows.
$ws.Range("A1:G32").AutoFilter() | Out-Null

$wb.Save()
